# Refresh the crypto price/volume table on Sheet1 with the latest
# scraped values (GitHub Actions data refresh, 2023-03-09).
# Price-like text cells that resemble plain numbers are written with a
# leading apostrophe (quote-prefix) so Excel keeps them as text, exactly
# as the source data has them (e.g. "1.001", "290.24"); the style is
# then reset to Normal so no stray number-format is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '21.687.12'
$ws.Range('E2').Value = '  -2.02%  '
$ws.Range('D3').Value = '1.536.76'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').Value = "'1.001"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'290.24"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('D8').Value = "'0.3195"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').Value = "'43.40"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('D10').Value = "'0.07200"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.08%  '
$ws.Range('D11').Value = "'1.075"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('D12').Value = "'1.001"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = "'5.760"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = "'18.43"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.81%  '
$ws.Range('D15').Value = "'6.627"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '1.534.71'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').Value = "'0.00001095"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.28%  '
$ws.Range('D18').Value = "'0.06617"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = "'84.16"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = "'0.9970"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = "'6.150"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('D22').Value = "'15.54"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('D23').Value = "'10.82"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.61%  '
$ws.Range('D24').Value = "'2.370"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').Value = '21.684.19'
$ws.Range('E25').Value = '  -2.07%  '
$ws.Range('D26').Value = "'2.387"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.62%  '
$ws.Range('D27').Value = "'150.87"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('D28').Value = "'18.50"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('D29').Value = "'4.892"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').Value = '1.707.00'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = "'117.61"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('D32').Value = "'6.075"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.56%  '
$ws.Range('D33').Value = "'0.9653"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -9.16%  '
$ws.Range('D34').Value = "'0.08098"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.50%  '
$ws.Range('D35').Value = "'5.197"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').Value = "'8.499"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.99%  '
$ws.Range('D37').Value = "'1.497"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.82%  '
$ws.Range('D38').Value = "'0.02222"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('D39').Value = "'0.05964"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.67%  '
$ws.Range('E40').Value = '  +4.38%  '
$ws.Range('D41').Value = "'0.2048"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').Value = "'1.184"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('D43').Value = "'1.0000"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = "'0.5817"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.17%  '
$ws.Range('D45').Value = "'13.21"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').Value = "'3.728"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').Value = "'0.5585"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'1.891"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = "'1.158"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').Value = "'115.94"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.57%  '
$ws.Range('D51').Value = "'0.06718"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.21%  '
